$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 100, shifting the existing rows 100-128 down to 102-130.
$ws.Rows("100:101").Insert()

# Row 100: new weekly entry (Magnum)
$ws.Range("A100").Value = 2
$ws.Range("B100").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C100").Value = "Coquimbo"
$ws.Range("D100").Value = 44559
$ws.Range("E100").Value = 4
$ws.Range("F100").Value = 100112031
$ws.Range("G100").Value = "Poroto verde"
$ws.Range("H100").Value = "Magnum"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 640
$ws.Range("K100").Value = 22000
$ws.Range("L100").Value = 23000
$ws.Range("M100").Value = 22500
$ws.Range("N100").Value = "$/malla 25 kilos"
$ws.Range("O100").Value = "Provincia de Limarí"
$ws.Range("P100").Value = 900
$ws.Range("Q100").Value = 25
$ws.Range("R100").Value = "Hortaliza"

# Row 101: new weekly entry (Sin especificar)
$ws.Range("A101").Value = 2
$ws.Range("B101").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 44559
$ws.Range("E101").Value = 4
$ws.Range("F101").Value = 100112031
$ws.Range("G101").Value = "Poroto verde"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 500
$ws.Range("K101").Value = 28000
$ws.Range("L101").Value = 30000
$ws.Range("M101").Value = 29000
$ws.Range("N101").Value = "$/malla 25 kilos"
$ws.Range("O101").Value = "Provincia de Limarí"
$ws.Range("P101").Value = 1160
$ws.Range("Q101").Value = 25
$ws.Range("R101").Value = "Hortaliza"

# Make sure the date cells keep the date number format (style index 2 in the original file).
$ws.Range("D100").NumberFormat = $ws.Range("D102").NumberFormat
$ws.Range("D101").NumberFormat = $ws.Range("D102").NumberFormat
